$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$c1 = $ws.Range("B1")
$c1.Borders.Weight = 2          # xlThin
$c1.Font.Bold = $true
$c1.HorizontalAlignment = -4108 # xlCenter
$c1.VerticalAlignment = -4160   # xlTop

$c1.Copy()
$c2 = $ws.Range("A2")
$c2.PasteSpecial(-4122)         # xlPasteFormats
